$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 13 entirely (it held only the professor's name in B13/C13,
#    with no label in column A). Everything below shifts up by one row.
$ws.Rows("13:13").Delete()

# 2. Overwrite specific cells with their new (shuffled) content.
#    B10:C10 (under "Objetivos:") now shows the professor identification.
$ws.Range("B10").Value = "5840705 - Maria Lúcia Caetano Pinto da Silva"
$ws.Range("C10").Value = "5840705 - Maria Lúcia Caetano Pinto da Silva"

#    B13:C13 (under "Programa resumido:") now shows "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

#    B15:C15 (under "Programa:") now shows the activation date "01/01/2022".
#    Copy the value from B8/C8 (which already holds that exact text) using
#    PasteSpecial-values so Excel keeps it as text instead of auto-parsing
#    it into a date serial number, and reuses the existing shared string.
$ws.Range("B8:C8").Copy()
$ws.Range("B15:C15").PasteSpecial(-4163)
$excel.CutCopyMode = $false

#    B18:C18 (under "Método:") now shows the professor identification again.
$ws.Range("B18").Value = "5840705 - Maria Lúcia Caetano Pinto da Silva"
$ws.Range("C18").Value = "5840705 - Maria Lúcia Caetano Pinto da Silva"

#    B19:C19 (under "Critério:") now shows the old "Método" text.
$ws.Range("B19").Value = "Serão oferecidas aulas expositivas e práticas."
$ws.Range("C19").Value = "Serão oferecidas aulas expositivas e práticas."

#    B20:C20 (under "Norma de recuperação:") now shows the old "Critério" text.
$ws.Range("B20").Value = "Serão aplicadas duas provas escritas. Trabalhos em sala de aula, seminários e relatórios, poderão, a critério do docente, ser considerados como parte da nota da prova escrita."
$ws.Range("C20").Value = "Serão aplicadas duas provas escritas. Trabalhos em sala de aula, seminários e relatórios, poderão, a critério do docente, ser considerados como parte da nota da prova escrita."

#    B21:C21 (under "Bibliografia:") now shows the old "Norma de recuperação" text,
#    and the long bibliography text that used to be here is gone.
$ws.Range("B21").Value = "Será realizada uma prova escrita envolvendo o conteúdo do semestre todo."
$ws.Range("C21").Value = "Será realizada uma prova escrita envolvendo o conteúdo do semestre todo."
